$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was inserted at row 208 (Fecha 2022-08-04 / serial
# 44777), pushing every existing record from the old row 208 down through the
# old last row (334) down by one row (209-335). Insert() shifts the existing
# rows down and copies the row's formatting (keeps the date style on column D).
$ws.Rows(208).Insert()

$ws.Cells.Item(208, 1).Value = 5
$ws.Cells.Item(208, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(208, 3).Value = "Maule"
$ws.Cells.Item(208, 4).Value = 44777
$ws.Cells.Item(208, 5).Value = 7
$ws.Cells.Item(208, 6).Value = 100112003
$ws.Cells.Item(208, 7).Value = "Ajo"
$ws.Cells.Item(208, 8).Value = "Chino"
$ws.Cells.Item(208, 9).Value = "Primera"
$ws.Cells.Item(208, 10).Value = 200
$ws.Cells.Item(208, 11).Value = 28000
$ws.Cells.Item(208, 12).Value = 28000
$ws.Cells.Item(208, 13).Value = 28000
$ws.Cells.Item(208, 14).Value = "$/malla 10 kilos"
$ws.Cells.Item(208, 15).Value = "China"
$ws.Cells.Item(208, 16).Value = 2800
$ws.Cells.Item(208, 17).Value = 10
$ws.Cells.Item(208, 18).Value = "Hortaliza"
